$d = $word.ActiveDocument

# Locate the end of the existing sentence " ... می‌باشد؟" (right before the
# _GoBack bookmark) so the new sentence can be appended immediately after it,
# inside the same paragraph.
$found = $d.Content.Find.Execute("می‌باشد؟", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertAt = $d.Content.End

$target = $d.Content
$target.Find.Execute("می‌باشد؟") | Out-Null
$e = $target.End

# Build the new sentence as three runs with the exact same run-level
# formatting Word itself would produce when typing Persian text, an English
# term and more Persian text back to back inside an RTL paragraph:
#   1) Persian lead-in  -> complex-script hint + rtl
#   2) "foreign key"    -> plain (no rtl / no cs hint), still bidi lang tag
#   3) Persian tail     -> complex-script hint + rtl
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> آیا می‌توان از </w:t></w:r>
<w:r><w:rPr><w:lang w:bidi="fa-IR"/></w:rPr><w:t>foreign key</w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> استفاده کرد؟</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

# InsertXML always creates a brand-new paragraph for the fragment, so insert
# it right after the found text (which pushes it past the zero-width
# _GoBack bookmark into its own paragraph) ...
$insertionPoint = $d.Range($e, $e)
$insertionPoint.InsertXML($xml)

# ... then delete the paragraph mark that now separates the original
# paragraph from the freshly inserted one so the three new runs become part
# of the same paragraph, right after the bookmark, exactly like the diff.
$paraMark = $d.Range($e, $e + 1)
$paraMark.Delete()
